$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AgeGroups")

# Row 26 (Masters M, group with I26..R26)
$ws.Range("I26").Value = ""
$ws.Range("J26").Value = "M45 49"
$ws.Range("K26").Value = "M49 55"
$ws.Range("L26").Value = "M55 59"
$ws.Range("M26").Value = "M61 65"
$ws.Range("N26").Value = "M67 71"
$ws.Range("O26").Value = "M73 79"
$ws.Range("P26").Value = "M81 87"
$ws.Range("Q26").Value = "M999 91"
$ws.Range("R26").Clear()

# Row 4 (Masters F)
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "F45 45"
$ws.Range("K4").Value = "F49 49"
$ws.Range("L4").Value = "F55 55"
$ws.Range("M4").Value = "F59 59"
$ws.Range("N4").Value = "F64 64"
$ws.Range("O4").Value = "F71 71"
$ws.Range("P4").Value = "F76 76"
$ws.Range("Q4").Value = "F999 80"

# Row 5
$ws.Range("I5").Value = ""

# Leave the selection where the author ended editing (Q4)
$ws.Activate()
$ws.Range("Q4").Select()
